# This workbook contains a weekly price log for "Bruselas (repollito)" at
# Vega Central Mapocho de Santiago. A new week's record needs to be added.
# In the source data this shows up as a new row being inserted right above
# the existing row 97 (which held date 45106 / 2023-06-29): the new row
# keeps all of the same J:R values as the old row 97 (same market/origin/
# price bracket) but carries a newer date (45132 / 2023-07-25). Every row
# from the old 97 through 104 shifts down by one (to 98 through 105).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above row 97; this automatically pushes former rows
# 97-104 down to 98-105 and extends the used range to A1:R105.
$ws.Rows.Item(97).Insert()

# Duplicate the row that is now at 98 (the original row 97 data) into the
# newly-created blank row 97.
$ws.Range("A97:R97").Value2 = $ws.Range("A98:R98").Value2

# The new row represents a later reporting date; update just the date cell.
$ws.Range("D97").Value2 = 45132
